# Weekly update: a new price observation is inserted at row 254 (pushing the
# existing rows 254-266 down to 255-267), extending the table by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 254, shifting rows 254:266 down to 255:267.
$ws.Rows("254").Insert()

# Populate the newly inserted row 254 with the new observation's data.
$ws.Range("A254").Value = 7
$ws.Range("B254").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C254").Value = 'Ñuble'
$ws.Range("D254").Value = 45147
$ws.Range("E254").Value = 16
$ws.Range("F254").Value = 'Fruta'
$ws.Range("G254").Value = 100109
$ws.Range("H254").Value = 'Uva'
$ws.Range("I254").Value = 100109001
$ws.Range("J254").Value = 'Uva'
$ws.Range("K254").Value = 'Crimpson Seedless'
$ws.Range("L254").Value = 'Primera'
$ws.Range("M254").Value = 60
$ws.Range("N254").Value = 12000
$ws.Range("O254").Value = 12000
$ws.Range("P254").Value = 12000
$ws.Range("Q254").Value = '$/bandeja 8 kilos'
$ws.Range("R254").Value = "Región de O'Higgins"
$ws.Range("S254").Value = 1500
$ws.Range("T254").Value = 8
